$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.739.71'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.80%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.049.17'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.63%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.05'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.93%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.64'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.44%  '

$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.046.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.58%  '

$ws.Range("E9").Value = '  +3.90%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.153'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.32%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.20'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -12.58%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.478'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.37%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000231'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.62%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.10'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.05%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.546.73'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.48%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.767.46'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.66%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.050.35'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.48%  '

$ws.Range("E18").Value = '  +0.32%  '

$ws.Range("E19").Value = '  -1.38%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '487.19'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.14%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.11'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.21%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.681'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.38'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.72%  '

$ws.Range("E24").Value = '  -0.28%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.43'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.54%  '

$ws.Range("E26").Value = '  +0.07%  '

$ws.Range("E27").Value = '  +0.01%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.09'
$ws.Range("D28").Style = "Normal"

$ws.Range("E29").Value = '  -1.86%  '

$ws.Range("E30").Value = '  +0.16%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.23'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.41%  '

$ws.Range("E32").Value = '  -0.33%  '

$ws.Range("E33").Value = '  -1.73%  '

$ws.Range("E34").Value = '  -0.28%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.18'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.02%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.16'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.05%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0407'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.64%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '439.64'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.88%  '

$ws.Range("E39").Value = '  -2.07%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.016.34'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.10%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.75'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.25%  '

$ws.Range("E42").Value = '  +0.12%  '

$ws.Range("E43").Value = '  -1.23%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.270'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.66%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '27.58'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.60%  '

$ws.Range("E46").Value = '  +3.75%  '

$ws.Range("E48").Value = '  +0.20%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '117.76'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.31%  '

$ws.Range("E50").Value = '  -2.05%  '

$ws.Range("E51").Value = '  +0.30%  '
